# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.633.20'
$ws.Range('E2').Value = '  -3.57%  '
$ws.Range('D3').Value = '3.300.10'
$ws.Range('E3').Value = '  -5.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.17%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = '3.300.46'
$ws.Range('E8').Value = '  -5.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.569'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.67%  '
$ws.Range('E10').Value = '  -5.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.56'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.521'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.84'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -11.35%  '
$ws.Range('E14').Value = '  -7.85%  '
$ws.Range('D15').Value = '3.844.03'
$ws.Range('E15').Value = '  -5.55%  '
$ws.Range('D16').Value = '67.700.50'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.317.30'
$ws.Range('E17').Value = '  -5.23%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '545.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -10.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.116'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.788'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.08'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.46%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -11.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '30.50'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.40%  '
$ws.Range('B31').Value = 'Mantle'
$ws.Range('C31').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.29%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -11.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '577.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.05'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -11.29%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0441'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0894'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.71'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.135'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.77'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -21.66%  '
$ws.Range('D43').Value = '3.006.97'
$ws.Range('E43').Value = '  -10.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.273'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -10.79%  '
$ws.Range('D46').Value = '0.0₃0613'
$ws.Range('E46').Value = '  -16.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -13.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.119'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.65%  '
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '127.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.96%  '
